$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 257 (pushes the existing 257:284 block down to 260:287)
$ws.Rows("257:259").Insert()

# New week's data (fecha 45124 = 2023-07-17), same structure as the rest of the
# Femacal de La Calera / Chirimoya block.
$rows = @(
    @{ Row=257; L="Especial"; M=56; N=33000; O=33000; P=33000; S=3300 },
    @{ Row=258; L="Primera";  M=57; N=30000; O=30000; P=30000; S=3000 },
    @{ Row=259; L="Segunda";  M=50; N=27000; O=27000; P=27000; S=2700 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 3
    $ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 45124
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 5
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = '$/bandeja 10 kilos'
    $ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 10
}

Write-Output "Done"
